$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Rank"
$ws.Range("B1").Value = "City Name"
$ws.Range("C1").Value = "Overnight International Visitors (Millions)"
$ws.Range("D1").Value = "Year"

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 4).Value = 2010
}

$ws.Range("G10").Select()
